$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

$ws.Range("B1").NumberFormat = "@"
$ws.Range("B1").Value = "00060992"

$ws.Range("F1").Value = "Hector Venzor"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "300006848"

$ws.Range("F3").Value = "Hector Venzor"

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "0880011962"

$ws.Range("E10").Value = "said safe to fly i swiped for it to fly went up and shot it self backwards into a fence propellers broke and crack on side of a leg on body"
